$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the H7 cell (Midnight Solution row / Dragonrot column) to the new value "Granite"
$ws.Range("H7").Value = "Granite"

# Update the selected cell to match the recorded sheet view selection
$ws.Range("H8").Select()
